# The document contains three occurrences of a pattern like:
#   <id>p163v_1</id>
# each one split across three separate runs (one run for the literal
# "<id>" text, one run for the "p163v_N" id value, one run for the
# literal "</id>" text). The edit collapses each of those three-run
# groups into a single run whose text is the full "<id>p163v_N</id>"
# string, carrying the formatting of the original first/third run
# (Courier New, color 7f6000, sz/szCs 18).
#
# Word's Find treats text across run boundaries as one continuous
# string, so Find.Execute locates the full "<id>p163v_N</id>" text even
# though it is split across runs. Re-assigning Range.Text then replaces
# the whole found range with a single run using the formatting of the
# range's first character - exactly the desired result. Setting the
# text to a throwaway placeholder first (in case the target text is
# already a single contiguous run in some edge case) guarantees the
# follow-up assignment is treated as a genuine change and actually
# performs the run-collapsing rewrite.

$d = $word.ActiveDocument

$ids = @("p163v_1", "p163v_2", "p163v_3")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"

    $rng = $d.Content
    $rng.Start = 0
    $rng.End = $d.Content.End

    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($found) {
        $rng.Text = "__tmp_placeholder__$id__"
        $rng.Text = $old
    }
}
